# Update "Data Test" sign-in test data for web, mobile, and api
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove only the hyperlink currently on C5 (archen@gmail.com) ---
# (its relationship slot will be reused once a new hyperlink is added below,
# which is what naturally renumbers the remaining C3 hyperlink's r:id)
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$5') {
        $hl.Delete()
    }
}

# --- Update the test data values ---
$ws.Range("C3").Value = "joongarch@gmail.com"
$ws.Range("C4").Value = "joong"
$ws.Range("D4").Value = "Joong123!"
$ws.Range("C5").Value = "joong@gmail.com"
$ws.Range("D5").Value = "Joong13!"

# --- Re-create the hyperlink on C5 ---
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:joong@gmail.com")

# Restore C5's original "Hyperlink" cell style (Hyperlinks.Add applies its
# own default formatting, which would otherwise introduce a new style).
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column C got a bit wider to fit the new email addresses ---
$ws.Columns.Item(3).ColumnWidth = 21.42

# --- Selection / view state at save time ---
$ws.Range("E11").Select() | Out-Null
